# Word COM-interop script: merges "  " + "Can be specified with " runs into
# a single run in 9 table cells, appends a new "Can be specified with
# ${Property} notation." sentence to the Alias cell, and relocates the
# "_GoBack" bookmark pair (and the "replaceValue" bookmark end) from the
# trailing empty paragraph after the second table into the tables proper.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Collapse the "  " / "Can be specified with " run pair into a single
#    run (same visible text) in each of the 9 affected table cells.
#    (The very first "double quotes... special characters." cell already
#    has these merged and must stay untouched.)
# ---------------------------------------------------------------------

$anchors = @(
  "The A part (basin name) to match, using * as a wildcard.  The location type part of the TSTool time series identifier is set to this value.  Can be specified with ",
  "The B part (location) to match, using * as a wildcard.  The location identifier part of the TSTool time series identifier is set to this value.  Can be specified with ",
  "The C part (parameter) to match, using * as a wildcard.  The TSTool data type is set to this value.  Can be specified with ",
  "The E part (interval) to match, using * as a wildcard.  Can be specified with ",
  "The F part (scenario) to match, using * as a wildcard.  Can be specified with ",
  "If specified, this will be used instead of the A-F parameters.  Can be specified with ",
  "Starting date/time to read data, in precision consistent with data.  Can be specified with ",
  "Ending date/time to read data, in precision consistent with data.  Can be specified with ",
  "This is useful when only Bpart is desired as the location identifier.  Can be specified with "
)

foreach ($anchor in $anchors) {
    $rng = $d.Content
    $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, $anchor, 2)
}

# ---------------------------------------------------------------------
# 2) Append the new "Can be specified with ${Property} notation." sentence
#    after " is scenario." in the Alias row's description cell.
# ---------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute(" is scenario.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("  ")
$rng.Collapse(0)
$rng.InsertAfter("Can be specified with ")
$rng.Collapse(0)
$rng.InsertAfter('${Property}')
$rng.Style = "RTiSWDocLiteralText"
$rng.Collapse(0)
$rng.InsertAfter(" notation.")

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark: it now starts/ends right after the new
#    sentence we just inserted (instead of in the trailing paragraph).
# ---------------------------------------------------------------------

$d.Bookmarks("_GoBack").Delete()
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)

# ---------------------------------------------------------------------
# 4) Move the end of the "replaceValue" bookmark so it closes at the end
#    of the second table (sample command table) instead of the trailing
#    paragraph after it.
# ---------------------------------------------------------------------

$tbl2 = $d.Tables(2)
$endRng = $tbl2.Range
$endRng.Collapse(0)

$startBm = $d.Bookmarks("replaceValue").Range
$startRng = $d.Range($startBm.Start, $startBm.Start)

$d.Bookmarks("replaceValue").Delete()
$d.Bookmarks.Add("replaceValue", $d.Range($startRng.Start, $endRng.Start))
